# CloseAllApplications.xaml and InitAllApplications.xaml no longer contain a
# workblock by default (main level + service levels). Workflows that still
# contain a workblock no longer pass in_wbType, since the workblock name is
# now decided inside the workflow. Reflect this in the Config.xlsx
# "Workblocks" sheet: keep only Init / GetTransactionData / ProcessTransaction
# workblock rows, clear out the rest, and update the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Row 5-6 used to describe the "CloseAppsRecover" workblock; it now
# describes "GetTransactionData" (shifted up from rows 7-8).
$ws.Range("A5").Value = "wbGetTransactionData_Type"
$ws.Range("B5").Value = "GetData"
$ws.Range("A6").Value = "wbGetTransactionData_SuppressSuccessful"

# Row 7-8 used to describe "GetTransactionData"; it now describes
# "ProcessTransaction" (shifted up from rows 9-10).
$ws.Range("A7").Value = "wbProcessTransaction_Type"
$ws.Range("B7").Value = "Process"
$ws.Range("A8").Value = "wbProcessTransaction_SuppressSuccessful"

# Everything from row 9 down (old ProcessTransaction, CloseAllApplications,
# InitAllApplications, Process rows) is no longer needed - clear it out.
$ws.Range("A9:C16").ClearContents()

# Bring Workblocks to the front (matches the new activeTab / tabSelected /
# selection state captured in the diff).
$ws.Activate() | Out-Null
$ws.Range("A3:C8").Select() | Out-Null
